$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(124, 1).Value = "20/11/2021"
$ws.Cells.Item(124, 2).Value = "HOLIDAY"

$ws.Cells.Item(125, 1).Value = "21/11/2021"
$ws.Cells.Item(125, 2).Value = "HOLIDAY"

$ws.Cells.Item(126, 1).Value = "22/11/2021"
$ws.Cells.Item(126, 2).Value = "Internal discussion on Memory management concepts"
$ws.Cells.Item(126, 3).Value = "Preparing ppt on Memory Management"

$ws.Cells.Item(127, 2).Value = "Prepared to give session on OpenMax"
$ws.Cells.Item(127, 3).Value = "Revising on basic C-DS-OS concepts"

$ws.Range("C127").Select()
